# Fill in the first sale line (row 7/8) with a product sold on
# 2025-08-23, and bump the footer timestamp from 9:25 AM to 9:26 AM.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: sale line item ------------------------------------------------

# A7: plain counter number
$ws.Range("A7").Value = 1

# C7:G7 and N7:O7 share one style (General number format) that switches to
# Text format once real data lands in it; apply to every cell of that style
# (as two contiguous blocks) so the workbook keeps a single shared style.
$ws.Range("C7:G7").NumberFormat = "@"
$ws.Range("N7:O7").NumberFormat = "@"
$ws.Range("C7").Value = "PROPAMETHONE TOP. CREAM. 20 GM"
$ws.Range("N7").Value = "30.00"

# H7:K7 share another style (General -> Text).
$ws.Range("H7:K7").NumberFormat = "@"
$ws.Range("H7").Value = "1:0"

# L7:M7 keep their original (numeric) number format but now hold a literal
# text value, so flip to Text only long enough to store it, then restore.
$origFmtL7 = $ws.Range("L7").NumberFormat
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = "1"
$ws.Range("L7").NumberFormat = $origFmtL7

# P7 likewise keeps its original numeric display format.
$origFmtP7 = $ws.Range("P7").NumberFormat
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "30.0000"
$ws.Range("P7").NumberFormat = $origFmtP7

# Q7 switches its style (General -> Text) permanently, reusing the same
# text value as H7.
$ws.Range("Q7").NumberFormat = "@"
$ws.Range("Q7").Value = "1:0"

# --- Row 8: total ----------------------------------------------------------
$ws.Range("P8").Value = 30

# --- Footer timestamp --------------------------------------------------------
$ws.Range("A9").Value = "Saturday, 23 August, 2025 9:26 AM"
